$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so numeric-looking
# strings (e.g. "1.003", "28.221.76") are not coerced into numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '28.221.76'
$ws.Range('E2').Value = '  +3.03%  '
$ws.Range('D3').Value = '1.810.46'
$ws.Range('E3').Value = '  +4.07%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '327.48'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.4343'
$ws.Range('E7').Value = '  +2.55%  '
$ws.Range('D8').Value = '0.3651'
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('D9').Value = '44.85'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').Value = '0.07652'
$ws.Range('E10').Value = '  +3.23%  '
$ws.Range('D11').Value = '1.140'
$ws.Range('E11').Value = '  +2.72%  '
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '21.99'
$ws.Range('E13').Value = '  +2.93%  '
$ws.Range('D14').Value = '6.279'
$ws.Range('E14').Value = '  +3.17%  '
$ws.Range('D15').Value = '7.502'
$ws.Range('E15').Value = '  +4.59%  '
$ws.Range('D16').Value = '1.826.53'
$ws.Range('E16').Value = '  +4.76%  '
$ws.Range('D17').Value = '93.59'
$ws.Range('E17').Value = '  +7.33%  '
$ws.Range('D18').Value = '0.00001079'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('D19').Value = '0.06564'
$ws.Range('E19').Value = '  +6.20%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = '17.42'
$ws.Range('E21').Value = '  +3.52%  '
$ws.Range('D22').Value = '6.241'
$ws.Range('E22').Value = '  +2.41%  '
$ws.Range('D23').Value = '28.242.22'
$ws.Range('E23').Value = '  +2.97%  '
$ws.Range('D24').Value = '11.55'
$ws.Range('E24').Value = '  -0.55%  '
$ws.Range('D25').Value = '2.077'
$ws.Range('E25').Value = '  -10.59%  '
$ws.Range('D26').Value = '161.58'
$ws.Range('E26').Value = '  +6.47%  '
$ws.Range('D27').Value = '20.61'
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('D28').Value = '2.035.83'
$ws.Range('E28').Value = '  +4.97%  '
$ws.Range('D29').Value = '2.274'
$ws.Range('E29').Value = '  -3.31%  '
$ws.Range('E30').Value = '  +1.83%  '
$ws.Range('D31').Value = '1.203'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').Value = '5.912'
$ws.Range('E32').Value = '  +4.35%  '
$ws.Range('D33').Value = '0.09150'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').Value = '3.454'
$ws.Range('E34').Value = '  -6.09%  '
$ws.Range('D35').Value = '12.96'
$ws.Range('E35').Value = '  +2.64%  '
$ws.Range('D36').Value = '0.02346'
$ws.Range('E36').Value = '  +2.78%  '
$ws.Range('D37').Value = '0.2168'
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('D38').Value = '5.173'
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('D39').Value = '0.6553'
$ws.Range('E39').Value = '  +2.53%  '
$ws.Range('D40').Value = '0.06177'
$ws.Range('E40').Value = '  +1.87%  '
$ws.Range('D41').Value = '1.191'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').Value = '8.088'
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('D43').Value = '1.437'
$ws.Range('E43').Value = '  +1.08%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '13.86'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('D46').Value = '0.6080'
$ws.Range('E46').Value = '  +3.67%  '
$ws.Range('D47').Value = '3.746'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('D48').Value = '125.43'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D49').Value = '2.009'
$ws.Range('E49').Value = '  +3.01%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.07006'
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '1.154'
$ws.Range('E51').Value = '  +2.81%  '
